# Remove the 4x3 numeric table that follows the final body paragraph
# (just before the section properties).
$d = $word.ActiveDocument

if ($d.Tables.Count -gt 0) {
    $d.Tables(1).Delete()
}
